$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 69, shifting existing rows 69:102 down to 70:103
$ws.Rows(69).Insert()

# Populate the newly inserted row 69 with the new data record
$ws.Range("A69").Value = 4
$ws.Range("B69").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C69").Value = "Los Lagos"
$ws.Range("D69").Value = 44466
$ws.Range("E69").Value = 10
$ws.Range("F69").Value = "Fruta"
$ws.Range("G69").Value = 100102
$ws.Range("H69").Value = "Cítricos"
$ws.Range("I69").Value = 100102004
$ws.Range("J69").Value = "Mandarina"
$ws.Range("K69").Value = "Murcott"
$ws.Range("L69").Value = "Primera"
$ws.Range("M69").Value = 400
$ws.Range("N69").Value = 6500
$ws.Range("O69").Value = 6500
$ws.Range("P69").Value = 6500
$ws.Range("Q69").Value = "$/bandeja 10 kilos"
$ws.Range("R69").Value = "Provincia de Limarí"
$ws.Range("S69").Value = 650
$ws.Range("T69").Value = 10
